$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'close for gym for woman',
    'cloth capris',
    'cloth frisbee',
    'cloth nighttime training pants',
    'cloth shorts for women',
    'cloth tights for women',
    'clothes capris',
    'clothes fitness',
    'clothes for adult women',
    'clothes for adults',
    'clothes for exercise',
    'clothes for gym women',
    'clothes for women',
    'clothes for women gym',
    'clothes for women over 50',
    'clothes for yoga women',
    'clothes line post',
    'clothes pants',
    'clothes tape to hold cloths in place',
    'clothes to wear on a cruise',
    'clothes women',
    'clothes women black',
    'clothes yoga pants',
    'clothing model stand',
    'clothing tape to keep clothes in place',
    'club dresses for women 2017',
    'club go generator',
    'club roots',
    'club tights for women',
    'club wear womens',
    'clubs leggings',
    'cold bike pants',
    'cold bike shirt',
    'cold compression gear',
    'cold compression knee',
    'cold compression leggings women',
    'cold compression pants',
    'cold compression shirt',
    'cold compression shirts for women',
    'cold compression therapy knee',
    'cold core wrap',
    'cold gear',
    'cold gear compression',
    'cold gear for men',
    'cold gear legging',
    'cold gear leggings',
    'cold gear leggings for women',
    'cold gear leggings women',
    'cold gear pants for women',
    'cold gear running tights women',
    'cold house warning light',
    'cold knee compression',
    'cold running gear',
    'cold running gear men',
    'cold running gear women',
    'cold sore',
    'cold weather biking gear',
    'cold weather capris women',
    'cold weather compression',
    'cold weather compression gear women',
    'cold weather compression pants men',
    'cold weather cycling gear',
    'cold weather essentials for women',
    'cold weather gear for biking',
    'cold weather gear for football',
    'cold weather gear running',
    'cold weather gear women',
    'cold weather high visibility work wear',
    'cold weather hiking pants',
    'cold weather hiking pants women',
    'cold weather hiking vest',
    'cold weather jogging gear for women',
    'cold weather leggings',
    'cold weather leggings for women',
    'cold weather leggings women running',
    'cold weather pants',
    'cold weather pants women',
    'cold weather riding pants',
    'cold weather run gear',
    'cold weather running',
    'cold weather running clothes',
    'cold weather running clothes for women',
    'cold weather running gear',
    'cold weather running gear for women',
    'cold weather running gear men',
    'cold weather running gear women',
    'cold weather running gear women top',
    'cold weather running gear womens',
    'cold weather running leggings',
    'cold weather running pants',
    'cold weather running pants with pockets for women',
    'cold weather running pants women',
    'cold weather running tights',
    'cold weather sports gear',
    'cold weather survival clothes',
    'cold weather tights',
    'cold weather tights for women',
    'cold weather women',
    'cold weather womens running gear',
    'cold weather womens running pants'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

